$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was M, now B)
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9444444444444444
$ws.Range("C2").Value = 0.9444444444444444
$ws.Range("D2").Value = 0.9444444444444444
$ws.Range("E2").Value = 36

# Row 3 (was B, now M)
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.9047619047619048
$ws.Range("C3").Value = 0.9047619047619048
$ws.Range("D3").Value = 0.9047619047619048
$ws.Range("E3").Value = 21

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.9298245614035088
$ws.Range("C4").Value = 0.9298245614035088
$ws.Range("D4").Value = 0.9298245614035088
$ws.Range("E4").Value = 0.9298245614035088

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9246031746031746
$ws.Range("C5").Value = 0.9246031746031746
$ws.Range("D5").Value = 0.9246031746031746

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9298245614035088
$ws.Range("C6").Value = 0.9298245614035088
$ws.Range("D6").Value = 0.9298245614035088
